# Update computed profit/price figures across the Brynhildr_Profits sheets
# (scheduled-runner refresh of market data: currentAveragePrice / Leve price / profit columns).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value2 = 400.66666
$ws.Range("J4").Value2 = 1002
$ws.Range("L4").Value2 = 1002
$ws.Range("N4").Value2 = -1230
$ws.Range("H5").Value2 = 208
$ws.Range("I5").Value2 = 130
$ws.Range("J5").Value2 = 325
$ws.Range("K5").Value2 = 130
$ws.Range("L5").Value2 = 325
$ws.Range("M5").Value2 = -15
$ws.Range("N5").Value2 = -555
$ws.Range("H62").Value2 = 12800
$ws.Range("J62").Value2 = 5000
$ws.Range("L62").Value2 = 5000
$ws.Range("N62").Value2 = -6248
$ws.Range("H65").Value2 = 12800
$ws.Range("J65").Value2 = 5000
$ws.Range("L65").Value2 = 25000
$ws.Range("N65").Value2 = -31240
$ws.Range("H69").Value2 = 7113.4
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 7113.4
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 21340.2
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value2 = -23088.2
$ws.Range("H72").Value2 = 7113.4
$ws.Range("I72").Value2 = 0
$ws.Range("J72").Value2 = 7113.4
$ws.Range("K72").Value2 = 0
$ws.Range("L72").Value2 = 64020.6
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value2 = -72756.60000000001
$ws.Range("H101").Value2 = 483.5
$ws.Range("I101").Value2 = 311.5
$ws.Range("K101").Value2 = 934.5
$ws.Range("M101").Value2 = 687.5
$ws.Range("H112").Value2 = 3380.889
$ws.Range("I112").Value2 = 0
$ws.Range("K112").Value2 = 0
$ws.Range("M112").ClearContents()
$ws.Range("H132").Value2 = 7348.9023
$ws.Range("I132").Value2 = 7594.212
$ws.Range("K132").Value2 = 22782.636
$ws.Range("M132").Value2 = -20252.636
$ws.Range("H138").Value2 = 3606.2646
$ws.Range("I138").Value2 = 3044
$ws.Range("J138").Value2 = 3752.037
$ws.Range("K138").Value2 = 9132
$ws.Range("L138").Value2 = 11256.111
$ws.Range("M138").Value2 = -3992
$ws.Range("N138").Value2 = -21536.111
$ws.Range("H141").Value2 = 3663.2903
$ws.Range("I141").Value2 = 1751.3684
$ws.Range("J141").Value2 = 6690.5
$ws.Range("K141").Value2 = 5254.1052
$ws.Range("L141").Value2 = 20071.5
$ws.Range("M141").Value2 = -74.10519999999997
$ws.Range("N141").Value2 = -30431.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H58").Value2 = 0
$ws.Range("I58").Value2 = 0
$ws.Range("K58").Value2 = 0
$ws.Range("M58").ClearContents()
$ws.Range("H61").Value2 = 2782237
$ws.Range("I61").Value2 = 4565.2
$ws.Range("J61").Value2 = 16670596
$ws.Range("K61").Value2 = 4565.2
$ws.Range("L61").Value2 = 16670596
$ws.Range("M61").Value2 = -4353.2
$ws.Range("N61").Value2 = -16671020
$ws.Range("H74").Value2 = 1329182.2
$ws.Range("I74").Value2 = 1686346.4
$ws.Range("K74").Value2 = 1686346.4
$ws.Range("M74").Value2 = -1685472.4
$ws.Range("H77").Value2 = 1329182.2
$ws.Range("I77").Value2 = 1686346.4
$ws.Range("K77").Value2 = 8431732
$ws.Range("M77").Value2 = -8427364
$ws.Range("H102").Value2 = 1020.1539
$ws.Range("I102").Value2 = 1020.1539
$ws.Range("K102").Value2 = 1020.1539
$ws.Range("M102").Value2 = 601.8461
$ws.Range("H110").Value2 = 947.7143
$ws.Range("I110").Value2 = 735.2381
$ws.Range("K110").Value2 = 735.2381
$ws.Range("M110").Value2 = 1309.7619
$ws.Range("H122").Value2 = 2663
$ws.Range("I122").Value2 = 2489
$ws.Range("K122").Value2 = 7467
$ws.Range("M122").Value2 = -5017
$ws.Range("H136").Value2 = 2782237
$ws.Range("I136").Value2 = 4565.2
$ws.Range("J136").Value2 = 16670596
$ws.Range("K136").Value2 = 13695.6
$ws.Range("L136").Value2 = 50011788
$ws.Range("M136").Value2 = -11145.6
$ws.Range("N136").Value2 = -50016888

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value2 = 7955.0625
$ws.Range("I64").Value2 = 12421.333
$ws.Range("K64").Value2 = 12421.333
$ws.Range("M64").Value2 = -12196.333
$ws.Range("H67").Value2 = 7955.0625
$ws.Range("I67").Value2 = 12421.333
$ws.Range("K67").Value2 = 12421.333
$ws.Range("M67").Value2 = -11641.333
$ws.Range("H94").Value2 = 215839.5
$ws.Range("I94").Value2 = 8712.9
$ws.Range("J94").Value2 = 1251472.5
$ws.Range("K94").Value2 = 8712.9
$ws.Range("L94").Value2 = 1251472.5
$ws.Range("M94").Value2 = -8261.9
$ws.Range("N94").Value2 = -1252374.5
$ws.Range("H107").Value2 = 849.62067
$ws.Range("I107").Value2 = 659.15
$ws.Range("K107").Value2 = 659.15
$ws.Range("M107").Value2 = 1260.85

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value2 = 19366.062
$ws.Range("I99").Value2 = 27139.75
$ws.Range("K99").Value2 = 27139.75
$ws.Range("M99").Value2 = -25641.75
$ws.Range("H103").Value2 = 10269.083
$ws.Range("I103").Value2 = 10269.083
$ws.Range("K103").Value2 = 10269.083
$ws.Range("M103").Value2 = -9097.083000000001
$ws.Range("H126").Value2 = 19366.062
$ws.Range("I126").Value2 = 27139.75
$ws.Range("K126").Value2 = 81419.25
$ws.Range("M126").Value2 = -78949.25
$ws.Range("H134").Value2 = 2277.4062
$ws.Range("I134").Value2 = 1906.5927
$ws.Range("J134").Value2 = 4279.8
$ws.Range("K134").Value2 = 5719.7781
$ws.Range("L134").Value2 = 12839.4
$ws.Range("M134").Value2 = -3184.7781
$ws.Range("N134").Value2 = -17909.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value2 = 55.93182
$ws.Range("J38").Value2 = 105.181816
$ws.Range("L38").Value2 = 315.545448
$ws.Range("N38").Value2 = -1009.545448
$ws.Range("H92").Value2 = 289
$ws.Range("J92").Value2 = 318.33334
$ws.Range("L92").Value2 = 955.0000200000001
$ws.Range("N92").Value2 = -3451.00002
$ws.Range("H132").Value2 = 4414.9355
$ws.Range("J132").Value2 = 6024.524
$ws.Range("L132").Value2 = 54220.716
$ws.Range("N132").Value2 = -59280.716

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value2 = 21993
$ws.Range("I18").Value2 = 21993
$ws.Range("K18").Value2 = 21993
$ws.Range("M18").Value2 = -21700
$ws.Range("H21").Value2 = 4893.25
$ws.Range("I21").Value2 = 4893.25
$ws.Range("K21").Value2 = 4893.25
$ws.Range("M21").Value2 = -4720.25
$ws.Range("H30").Value2 = 4893.25
$ws.Range("I30").Value2 = 4893.25
$ws.Range("K30").Value2 = 4893.25
$ws.Range("M30").Value2 = -4788.25
$ws.Range("H122").Value2 = 5181.6
$ws.Range("I122").Value2 = 5439
$ws.Range("J122").Value2 = 4867
$ws.Range("K122").Value2 = 16317
$ws.Range("L122").Value2 = 14601
$ws.Range("M122").Value2 = -13867
$ws.Range("N122").Value2 = -19501
$ws.Range("H132").Value2 = 10735.756
$ws.Range("I132").Value2 = 9674.675999999999
$ws.Range("J132").Value2 = 15643.25
$ws.Range("K132").Value2 = 29024.028
$ws.Range("L132").Value2 = 46929.75
$ws.Range("M132").Value2 = -26494.028
$ws.Range("N132").Value2 = -51989.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value2 = 4636.6665
$ws.Range("I7").Value2 = 0
$ws.Range("J7").Value2 = 4636.6665
$ws.Range("K7").Value2 = 0
$ws.Range("L7").Value2 = 4636.6665
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value2 = -4860.6665
$ws.Range("H23").Value2 = 9450
$ws.Range("I23").Value2 = 9450
$ws.Range("K23").Value2 = 9450
$ws.Range("M23").Value2 = -9220
$ws.Range("H46").Value2 = 3852.1
$ws.Range("I46").Value2 = 1039.8
$ws.Range("J46").Value2 = 4789.533
$ws.Range("K46").Value2 = 1039.8
$ws.Range("L46").Value2 = 4789.533
$ws.Range("M46").Value2 = -851.8
$ws.Range("N46").Value2 = -5165.533
$ws.Range("H68").Value2 = 1954.2222
$ws.Range("J68").Value2 = 2000
$ws.Range("L68").Value2 = 2000
$ws.Range("N68").Value2 = -3498
$ws.Range("H69").Value2 = 0
$ws.Range("J69").Value2 = 0
$ws.Range("L69").Value2 = 0
$ws.Range("N69").ClearContents()
$ws.Range("H71").Value2 = 1954.2222
$ws.Range("J71").Value2 = 2000
$ws.Range("L71").Value2 = 10000
$ws.Range("N71").Value2 = -17488
$ws.Range("H72").Value2 = 0
$ws.Range("J72").Value2 = 0
$ws.Range("L72").Value2 = 0
$ws.Range("N72").ClearContents()
$ws.Range("H122").Value2 = 3351.1667
$ws.Range("I122").Value2 = 3221.9
$ws.Range("J122").Value2 = 3997.5
$ws.Range("K122").Value2 = 9665.700000000001
$ws.Range("L122").Value2 = 11992.5
$ws.Range("M122").Value2 = -7215.700000000001
$ws.Range("N122").Value2 = -16892.5
$ws.Range("H126").Value2 = 4636.6665
$ws.Range("I126").Value2 = 0
$ws.Range("J126").Value2 = 4636.6665
$ws.Range("K126").Value2 = 0
$ws.Range("L126").Value2 = 13909.9995
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value2 = -18849.9995
$ws.Range("H132").Value2 = 4031923.2
$ws.Range("I132").Value2 = 5082708.5
$ws.Range("K132").Value2 = 15248125.5
$ws.Range("M132").Value2 = -15245595.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value2 = 12815532
$ws.Range("I136").Value2 = 2899540
$ws.Range("J136").Value2 = 50000504
$ws.Range("K136").Value2 = 8698620
$ws.Range("L136").Value2 = 150001512
$ws.Range("M136").Value2 = -8696070
$ws.Range("N136").Value2 = -150006612

Write-Output "Updated $($wb.Worksheets.Count) sheets with refreshed market data."
